{"js": "// Remove the \"Ver no Jupiter / Salvar em pdf / Salvar em docx\" paragraph,\n// the \"\u00a9 2020 ... Creative Commons Attribution\" paragraph, and the blank\n// paragraph that immediately followed them (the footer block that the\n// Jekyll site build no longer emits).\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\nconst items = paragraphs.items;\n\n// Locate the \"Ver no Jupiter ...\" paragraph by its distinctive text.\nlet startIndex = -1;\nfor (let i = 0; i < items.length; i++) {\n  if (items[i].text.indexOf(\"Ver no Jupiter\") !== -1) {\n    startIndex = i;\n    break;\n  }\n}\n\nif (startIndex !== -1) {\n  // Delete the three consecutive paragraphs: \"Ver no Jupiter ...\",\n  // \"\u00a9 2020 ...\", and the blank paragraph right after them. Deleting from\n  // the last one back to the first avoids disturbing the position of the\n  // paragraphs still queued for removal.\n  const endIndex = Math.min(startIndex + 2, items.length - 1);\n  for (let i = endIndex; i >= startIndex; i--) {\n    items[i].delete();\n  }\n  await context.sync();\n}\n", "ps1": "# Remove the \"Ver no Jupiter / Salvar em pdf / Salvar em docx\" paragraph,\n# the \"\u00a9 2020 ... Creative Commons Attribution\" paragraph, and the blank\n# paragraph right after them (the Jekyll-site footer block that the\n# rebuilt page no longer emits).\n$d = $word.ActiveDocument\n\n$searchRange = $d.Content\n$found = $searchRange.Find.Execute(\"Ver no Jupiter\")\n\nif ($found) {\n    $count = $d.Paragraphs.Count\n    $startIndex = -1\n    for ($i = 1; $i -le $count; $i++) {\n        $p = $d.Paragraphs.Item($i)\n        if ($p.Range.Start -le $searchRange.Start -and $p.Range.End -ge $searchRange.End) {\n            $startIndex = $i\n            break\n        }\n    }\n\n    if ($startIndex -ne -1) {\n        # Remove this paragraph plus the two that follow it (copyright line\n        # and the trailing blank paragraph).\n        $endIndex = [Math]::Min($startIndex + 2, $count)\n        $startPara = $d.Paragraphs.Item($startIndex)\n        $endPara = $d.Paragraphs.Item($endIndex)\n        $deleteRange = $d.Range($startPara.Range.Start, $endPara.Range.End)\n        $deleteRange.Delete()\n    }\n}\n"}
